$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: all Price (D) and Hora (G) values are stored as text in the sheet.
# A leading apostrophe is used so Excel keeps numeric-looking strings (e.g.
# "245.29", "8", "0.0001000") as text instead of coercing them to numbers,
# which matches the workbook's original inlineStr/text cell formatting.
$ws.Range("D2").Value = "'245.29"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'21.99"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'5.404"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.05973"
$ws.Range("G5").Value = "'8"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'6.399"
$ws.Range("G7").Value = "'8"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'0.9629"
$ws.Range("G9").Value = "'8"
$ws.Range("G10").Value = "'8"
$ws.Range("G11").Value = "'8"
$ws.Range("D12").Value = "'0.03415"
$ws.Range("G12").Value = "'8"
$ws.Range("D13").Value = "'0.03068"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.09423"
$ws.Range("G14").Value = "'8"
$ws.Range("D15").Value = "'4.003"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'0.001597"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'0.04825"
$ws.Range("G17").Value = "'8"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005942"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "'8"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006146"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("G19").Value = "'8"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005105"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Value = "'8"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009871"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("G21").Value = "'8"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001000"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("G22").Value = "'8"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.745"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("G23").Value = "'8"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.186"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").Value = "'8"
$ws.Range("G25").Value = "'8"
$ws.Range("G26").Value = "'8"
$ws.Range("D27").Value = "'0.0002462"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("G38").Value = "'8"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.03983"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.006383"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("G42").Value = "'8"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.005301"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.00005250"
$ws.Range("G45").Value = "'8"
$ws.Range("G46").Value = "'8"
$ws.Range("D47").Value = "'0.6702"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("G47").Value = "'8"
$ws.Range("D48").Value = "'0.02964"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "'8"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("G49").Value = "'8"
$ws.Range("G50").Value = "'8"
$ws.Range("G51").Value = "'8"
